$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: new "2.0" log entry --------------------------------------------

# Hour number
$ws.Range("A7").Value = 6

# "project versie" column (B): copy the formatting (right-aligned, General
# number format) from an existing entry first...
$ws.Range("B2").Copy()
$ws.Range("B7").PasteSpecial(-4122)

# ...then push in the literal text "2.0". A plain Value assignment of "2.0"
# gets auto-coerced to the number 2 by Excel, so instead build it with a
# text formula in a scratch cell and paste back just the resulting value -
# that keeps it a genuine text/shared-string cell.
$helper = $ws.Range("H1")
$helper.Formula = '="2"&"."&"0"'
$helper.Copy()
$ws.Range("B7").PasteSpecial(-4163)
$helper.ClearContents()

# "Beschrijving" column (C): copy formatting from an existing description
# cell (column C's default/autofill format is bold, so this avoids the new
# row picking that up), then set the real text.
$ws.Range("C2").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("C7").Value = "Dit uur heb ik als eerste een meme toegevoegd als je het project opent. Daarna heb ik een `"Pretty mode`" button gemaakt. Deze activeerd / deactiveerd de return yield null in de GenerateScript IEnumerator. Hierdoor is er of instant een maze gegenereerd, of visueel stapje voor stapje. als laatste heb ik een klein stukje code toegevoeg zodat de start positie van de camera gebaseerd is op de grote van de maze. Hierdoor zie je altijd de hele maze aan het begin."

# --- Row 8: blank spacer row, same alignment style as column B -------------
$ws.Range("B8").HorizontalAlignment = -4152

# --- View / selection changes ------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("C7").Select()
